$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.513.58"

$ws.Range("D3").Value = "1.869.17"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").Value = "'312.65"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").Value = "'0.4783"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("D8").Value = "'0.3781"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Range("E8").Value = "  +3.09%  "

$ws.Range("D9").Value = "'0.07363"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").Value = "'0.9365"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("D11").Value = "'20.76"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Range("E11").Value = "  +5.21%  "

$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("D13").Value = "1.877.33"
$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("E14").Value = "  +2.41%  "

$ws.Range("D15").Value = "'6.580"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "'90.75"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("D18").Value = "'0.000008899"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Range("E18").Value = "  +3.18%  "

$ws.Range("D19").Value = "'1.012"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "'14.94"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Range("E20").Value = "  +2.75%  "

$ws.Range("D21").Value = "27.515.52"

$ws.Range("D22").Value = "'5.136"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("D23").Value = "'10.73"
$ws.Cells.Item(23,4).Style = "Normal"

$ws.Range("D24").Value = "'1.956"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("D25").Value = "'153.94"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("E26").Value = "  +2.25%  "

$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").Value = "'116.00"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("D29").Value = "'5.003"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("D30").Value = "'0.08932"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("D31").Value = "'3.341"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("D32").Value = "'1.221"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Range("E32").Value = "  +4.01%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.619"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Range("E33").Value = "  +3.08%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7544"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").Value = "'2.691"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").Value = "'0.02052"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Range("E36").Value = "  +4.97%  "

$ws.Range("D37").Value = "'1.119"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.003"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05276"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").Value = "'0.5359"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "

$ws.Range("D41").Value = "'7.083"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "

$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("D43").Value = "'8.468"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Range("E43").Value = "  +3.14%  "

$ws.Range("D44").Value = "'10.66"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("D45").Value = "'0.4810"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D48").Value = "'102.81"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("D49").Value = "'67.49"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("D50").Value = "'0.06088"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "'0.9259"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Range("E51").Value = "  +4.66%  "
